$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 339, shifting existing rows 339:363 down to 340:364
$ws.Rows(339).Insert()

# Populate the newly inserted row 339 with the new record's data
$ws.Cells.Item(339, 1).Value = 10
$ws.Cells.Item(339, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(339, 3).Value = "La Araucanía"
$ws.Cells.Item(339, 4).Value = 44714
$ws.Cells.Item(339, 5).Value = 9
$ws.Cells.Item(339, 6).Value = "Fruta"
$ws.Cells.Item(339, 7).Value = 100108
$ws.Cells.Item(339, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(339, 9).Value = 100108002
$ws.Cells.Item(339, 10).Value = "Mango"
$ws.Cells.Item(339, 11).Value = "Sin especificar"
$ws.Cells.Item(339, 12).Value = "Primera"
$ws.Cells.Item(339, 13).Value = 900
$ws.Cells.Item(339, 14).Value = 11000
$ws.Cells.Item(339, 15).Value = 11000
$ws.Cells.Item(339, 16).Value = 11000
$ws.Cells.Item(339, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(339, 18).Value = "Brasil"
$ws.Cells.Item(339, 19).Value = 2750
$ws.Cells.Item(339, 20).Value = 4
